$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("杨瀚森"): add a new vocabulary-list row (row 10) ---------------
# Clear the stray truly-empty "t=s" placeholder cells (void index values left
# over from the source export) so they no longer emit empty <c> entries.
# I2 ("重新复习：") is real content and must stay untouched.
$ws1.Range("H2").ClearContents()
$ws1.Range("H3:I9").ClearContents()

# Copy row 9's formatting down into row 10, then overwrite with the new data.
$ws1.Range("A9:I9").Copy($ws1.Range("A10:I10"))
$ws1.Range("H10:I10").ClearContents()

$ws1.Cells.Item(10, 1).Value = 8
$ws1.Cells.Item(10, 2).Value = 43082
$ws1.Cells.Item(10, 3).Value = "高中单词9"
$ws1.Cells.Item(10, 4).Value = 0
$ws1.Cells.Item(10, 5).Value = 43082
$ws1.Cells.Item(10, 6).Value = 43083
$ws1.Cells.Item(10, 7).Value = 43085

$ws1.Columns.Item(5).ColumnWidth = 20.5

# --- Sheet2 ("尹嘉禾"): update the current-task log row (row 2) -------------
$ws2.Cells.Item(2, 1).Value = 0
$ws2.Cells.Item(2, 2).Value = 43082
$ws2.Cells.Item(2, 3).Value = "21天list1"
$ws2.Cells.Item(2, 4).Value = 0
$ws2.Cells.Item(2, 5).Value = 43082
$ws2.Cells.Item(2, 6).Value = 43083
$ws2.Cells.Item(2, 7).Value = 43085

$ws2.Range("E3").Select()

# Sheet1 is the active/visible tab in the saved workbook, so re-select it
# (and its last cell) after touching sheet2.
$ws1.Range("I9").Select()
